$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing hour values ---
$ws.Range("F18").Value = 0.75
$ws.Range("F20").Value = 2.5

# --- Pre-fill the numeric/date values for the two new tracked rows so the
#     SUM(F:F) formula in C5 picks them up before we touch formatting ---
$ws.Range("E21").Value = 43683
$ws.Range("F21").Value = 2
$ws.Range("E22").Value = 43684
$ws.Range("F22").Value = 4

# --- Copy the row-20 formatting down onto the two new rows ---
$ws.Range("E20:G20").Copy()
$ws.Range("E21:G22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 21 text ---
$ws.Range("G21").Value = "Arbeiten am ausgewählten Design (UI resizeable machen)"

# --- Row 22 text (multi-line, wrapped) ---
$ws.Range("G22").Value = "Probleme beim Arbeiten an der UI aufgetreten (UI ist aufgrund des Designs nicht so leicht resizeable zu machen,`nwäre mit höheren Arbeitsaufwand verbunden), Arbeiten an dieser Funktion auf späteren Zeitpunkt verschoben`nOptische Anpassungen an Design"
$ws.Range("G22").WrapText = $true
$ws.Range("G22").HorizontalAlignment = -4131  # xlLeft
$ws.Rows("22:22").RowHeight = 57.6

# --- Row 23 - final note, reuse the standard "Arbeit" column formatting ---
$ws.Range("G20").Copy()
$ws.Range("G23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("G23").Value = "Funktionalität der UI programmiert"

# --- Make column G a bit wider to accommodate the new wrapped text ---
$ws.Columns("G:G").ColumnWidth = 91

$excel.Calculate()

# --- Reflect the final selection like the author left it ---
[void]$ws.Range("G22").Select()

Write-Host "done"
